$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 144.783305
$ws.Range("H2").Value = 434.349915
$ws.Range("I2").Value = 0.2430046335191003
$ws.Range("J2").Value = 0.251012682214973
$ws.Range("Q2").Value = 12.03810441642833
$ws.Range("R2").Value = 108.342939747855
$ws.Range("S2").Value = 0.2430046335191003
$ws.Range("T2").Value = 0.251012682214973

# Row 3
$ws.Range("G3").Value = 82.24887099999999
$ws.Range("I3").Value = 0.1380466950572427
$ws.Range("J3").Value = 0.1425959278859072
$ws.Range("Q3").Value = 6.838637211875666
$ws.Range("S3").Value = 0.1380466950572427
$ws.Range("T3").Value = 0.1425959278859072

# Row 4
$ws.Range("G4").Value = 163.8590903333333
$ws.Range("H4").Value = 491.577271
$ws.Range("I4").Value = 0.2750214756820535
$ws.Range("J4").Value = 0.284084617144743
$ws.Range("Q4").Value = 13.62417330515856
$ws.Range("R4").Value = 122.617559746427
$ws.Range("S4").Value = 0.2750214756820535
$ws.Range("T4").Value = 0.284084617144743

# Row 5
$ws.Range("G5").Value = 57.0238095
$ws.Range("H5").Value = 114.047619
$ws.Range("I5").Value = 0.09570889357312636
$ws.Range("J5").Value = 0.06590860906562239
$ws.Range("Q5").Value = 4.741282656750499
$ws.Range("R5").Value = 28.447695940503
$ws.Range("S5").Value = 0.09570889357312636
$ws.Range("T5").Value = 0.06590860906562239

# Row 6
$ws.Range("G6").Value = 147.8896333333333
$ws.Range("H6").Value = 443.6689
$ws.Range("I6").Value = 0.2482183021684772
$ws.Range("J6").Value = 0.2563981636887546
$ws.Range("Q6").Value = 12.29638215658889
$ws.Range("R6").Value = 110.6674394093
$ws.Range("S6").Value = 0.2482183021684772
$ws.Range("T6").Value = 0.2563981636887546
